# "working on the new function of SHOCK"
#
# - Rename worksheet "A" -> "Z".
# - Make "Z" the active sheet/tab (was "Y"), with C10 selected
#   (and C1 scrolled into view as the pane's top-left cell).
# - Clear the leftover placeholder numbers (2,3,4,5) out of A3:A6 on "Z",
#   keeping their existing cell style.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("A")
$ws.Name = "Z"

# Drop the old index values from A3:A6 but keep the cell formatting.
$ws.Range("A3:A6").ClearContents()

# Activate "Z" and move the selection/view the way the commit left it.
$ws.Activate()
$ws.Range("C1").Select()
$ws.Range("C10").Select()

Write-Output "done"
